$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new "version 4.0" record -------------------------------------
# Order matters for shared-string table layout: the version text ("version
# 4.0") was registered before the longer commentary text in the source
# workbook, so write C12 before B12.
$ws.Range("C12").Value = "version 4.0"
$ws.Range("B12").Value = "Using mutex to avoid race condition and modularized code into different .cpp and .h files"
$ws.Range("D12").Value = 10000
$ws.Range("E12").Value = 9017
$ws.Range("F12").Formula = "=E12/1000"

# Match the formatting used by the row above it (same column styles as the
# rest of the table: B = wrap text, C = left aligned, D = number format).
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows(12).RowHeight = 27.6

# --- Tidy up C6:C7 formatting ----------------------------------------------
# These two cells carried a leftover "empty alignment" style with no actual
# effect; reset them back to the Normal style.
$ws.Range("C6:C7").Style = "Normal"

# --- Selection --------------------------------------------------------------
$ws.Range("B15").Select()
